# Updates the crypto price/volume table on Sheet1 (columns B-E, rows 2-51)
# to match the refreshed data pulled on Sun May  7 02:43:36 UTC 2023.
#
# Most rows only get new Price (D) / Volume 1h (E) figures. Rows 45-48
# also change Coin (B) and Link (C) because the coin list reshuffled:
#   PEPE/EnergySwap/RenderToken/Cronos -> EnergySwap/RenderToken/Cronos/PEPE
#
# Several Price values look numeric (e.g. "1.006", "19.10") but must be
# preserved as literal text (matching the trailing zeros / exact digits
# in the source feed), so those cells are switched to a text number
# format before the value is written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.027.32"
$ws.Range("E2").Value = "  -1.81%  "
$ws.Range("D3").Value = "1.905.20"
$ws.Range("E3").Value = "  -4.02%  "
$ws.Range("D4").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D5").Value = "325.12"
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("D6").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D7").Value = "0.4596"
$ws.Range("E7").Value = "  -1.26%  "
$ws.Range("D8").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D8").Value = "0.3808"
$ws.Range("E8").Value = "  -2.99%  "
$ws.Range("D9").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D9").Value = "0.07708"
$ws.Range("E9").Value = "  -2.75%  "
$ws.Range("D10").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D10").Value = "0.9733"
$ws.Range("E10").Value = "  -1.97%  "
$ws.Range("D11").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D11").Value = "21.92"
$ws.Range("E11").Value = "  -4.08%  "
$ws.Range("D12").Value = "1.944.21"
$ws.Range("E12").Value = "  -2.88%  "
$ws.Range("D13").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D13").Value = "6.932"
$ws.Range("E13").Value = "  -3.55%  "
$ws.Range("D14").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D14").Value = "5.652"
$ws.Range("E14").Value = "  -3.27%  "
$ws.Range("D15").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D15").Value = "0.07085"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D16").Value = "1.007"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D17").Value = "83.65"
$ws.Range("E17").Value = "  -4.48%  "
$ws.Range("D18").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D18").Value = "0.000009478"
$ws.Range("E18").Value = "  -4.62%  "
$ws.Range("D19").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D19").Value = "16.62"
$ws.Range("E19").Value = "  -4.03%  "
$ws.Range("D20").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").Value = "29.003.14"
$ws.Range("E21").Value = "  -2.00%  "
$ws.Range("D22").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D22").Value = "5.302"
$ws.Range("E22").Value = "  -4.41%  "
$ws.Range("D23").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D23").Value = "10.87"
$ws.Range("E23").Value = "  -2.86%  "
$ws.Range("D24").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D24").Value = "2.098"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D25").Value = "157.95"
$ws.Range("E25").Value = "  -0.60%  "
$ws.Range("D26").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D26").Value = "19.10"
$ws.Range("E26").Value = "  -2.58%  "
$ws.Range("D27").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D27").Value = "5.589"
$ws.Range("E27").Value = "  -3.85%  "
$ws.Range("D28").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D28").Value = "117.59"
$ws.Range("E28").Value = "  -1.70%  "
$ws.Range("D29").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D29").Value = "1.846"
$ws.Range("E29").Value = "  -2.87%  "
$ws.Range("D30").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D30").Value = "0.09257"
$ws.Range("E30").Value = "  -1.66%  "
$ws.Range("D31").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D31").Value = "0.8573"
$ws.Range("E31").Value = "  -4.48%  "
$ws.Range("D32").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D32").Value = "5.080"
$ws.Range("E32").Value = "  -2.71%  "
$ws.Range("D33").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D33").Value = "1.236"
$ws.Range("E33").Value = "  -7.31%  "
$ws.Range("D34").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D34").Value = "2.997"
$ws.Range("E34").Value = "  -6.14%  "
$ws.Range("D35").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D35").Value = "0.05672"
$ws.Range("E35").Value = "  -2.30%  "
$ws.Range("E36").Value = "  -3.05%  "
$ws.Range("D37").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D37").Value = "1.004"
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("E38").Value = "  -2.95%  "
$ws.Range("D39").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D39").Value = "0.5481"
$ws.Range("E39").Value = "  -4.30%  "
$ws.Range("D40").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D40").Value = "7.393"
$ws.Range("E40").Value = "  -5.70%  "
$ws.Range("D41").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D41").Value = "0.1751"
$ws.Range("E41").Value = "  -2.72%  "
$ws.Range("D42").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D42").Value = "9.300"
$ws.Range("E42").Value = "  -4.28%  "
$ws.Range("D43").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D43").Value = "2.765"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D44").Value = "0.5156"
$ws.Range("E44").Value = "  -3.76%  "
# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D45").Value = "11.22"
$ws.Range("E45").Value = "  -6.30%  "
# Row 46
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D46").Value = "2.082"
$ws.Range("E46").Value = "  -4.46%  "
# Row 47
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D47").Value = "0.06833"
$ws.Range("E47").Value = "  -1.61%  "
# Row 48
$ws.Range("B48").Value = "PEPE"
$ws.Range("C48").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D48").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D48").Value = "0.000002605"
$ws.Range("E48").Value = "  -21.08%  "
$ws.Range("D49").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D49").Value = "1.771"
$ws.Range("E49").Value = "  -2.86%  "
$ws.Range("D50").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D50").Value = "110.04"
$ws.Range("E50").Value = "  -3.66%  "
$ws.Range("D51").NumberFormat = "@"   # keep as text, not a number
$ws.Range("D51").Value = "1.003"
$ws.Range("E51").Value = "  +0.06%  "
